# Apply updated crypto price/volume figures to columns D (Price) and E (Volume 1h).
# A leading apostrophe forces Excel to store number-like strings (e.g. 1.001)
# as text, matching the original inlineStr cell type instead of being parsed as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.974.75"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.643.30"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'214.92"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'0.5214"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "'0.07677"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").Value = "1.646.66"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "'4.420"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "1.866.74"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "0.0₅8278"
$ws.Range("E16").Value = "  +3.55%  "
$ws.Range("D17").Value = "'64.71"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "25.986.01"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'4.703"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "'188.75"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "'10.17"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "'6.260"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'144.32"
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").Value = "'0.1222"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "'7.395"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'15.87"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").Value = "'0.05924"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("D31").Value = "'1.263"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'3.392"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "'3.398"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").Value = "'0.9929"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "'2.393"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'2.754"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'0.5635"
$ws.Range("E38").Value = "  -5.78%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").Value = "'5.869"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "1.030.52"
$ws.Range("D44").Value = "'98.83"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "1.792.92"
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'55.64"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'8.039"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'0.05144"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "'0.4220"
$ws.Range("E51").Value = "  -0.36%  "
